$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reports")

# New rows of report data to append, mirroring the existing layout
# (Request ID, Report Timestamp, Report)
$newRows = @(
    @{ Id = 26; Timestamp = "2025-04-26 21:48:13"; Report = "John Smith moved Its brand new from ford from No Location to shelf space 2.`nNow John Smith is Excited.`n" },
    @{ Id = 27; Timestamp = "2025-04-26 21:49:24"; Report = "John Smith took picture of Its brand new from ford.`nNow John Smith is Frustrated.`n" },
    @{ Id = 28; Timestamp = "2025-04-26 23:34:14"; Report = "John Smith added Nissan battery to the database.`nSerial Number is 12e12eknkndkfak.`nPart Number is 7.`nItem Type is 3.`nNow John Smith is Confident.`n" },
    @{ Id = 29; Timestamp = "2025-04-26 23:43:17"; Report = "John Smith added Suzuki Battery to the database.`nSerial Number is 78578assa87sa87as758.`nPart Number is 9.`nItem Type is 2.`nLocation is shelf space 1.`nNow John Smith is Bored.`n" }
)

$startRow = 27
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.Id
    $ws.Cells.Item($r, 2).Value = $row.Timestamp

    $cellC = $ws.Cells.Item($r, 3)
    $cellC.Value = $row.Report
    $cellC.WrapText = $true
}
